{"js": "// Apply benchmark-result updates to the single-column results table.\n// Each entry is [zero-based row index, new cell text]. Rows 43-45 (0-based,\n// i.e. the 44th-46th rows) previously held a tab-separated run of 10 values\n// each; they collapse down to just their first value, while rows 0-11\n// (0-based, the 1st-12th rows) receive updated figures.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"364\"],\n  [4, \"0.00003\"],\n  [5, \"0.00271\"],\n  [6, \"0.00019\"],\n  [7, \"0.00009\"],\n  [8, \"0.00027\"],\n  [9, \"0.00038\"],\n  [10, \"0.00045\"],\n  [11, \"0.08146\"],\n  [43, \"100\"],\n  [44, \"0.08\"],\n  [45, \"2439\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Apply benchmark-result updates to the single-column results table.\n# Each entry is [1-based row index, new cell text]. Rows 44-46 previously\n# held a tab-separated run of 10 values; they collapse down to just their\n# first value, while rows 1-12 receive updated figures.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @(1, \"0M\"),\n    @(2, \"0M\"),\n    @(3, \"0M\"),\n    @(4, \"364\"),\n    @(5, \"0.00003\"),\n    @(6, \"0.00271\"),\n    @(7, \"0.00019\"),\n    @(8, \"0.00009\"),\n    @(9, \"0.00027\"),\n    @(10, \"0.00038\"),\n    @(11, \"0.00045\"),\n    @(12, \"0.08146\"),\n    @(44, \"100\"),\n    @(45, \"0.08\"),\n    @(46, \"2439\")\n)\n\nforeach ($pair in $updates) {\n    $rowIndex = $pair[0]\n    $newText = $pair[1]\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $newText\n}\n"}
